$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in columns D/E hold numeric-looking text (prices, percentages) that
# must stay TEXT, matching the source inline-string cells. Force text typing
# by setting the number format to Text ('@') before assigning the value, then
# restore the cell style to 'Normal' so no stray formatting diff is left behind.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '307.41'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '-1.09%'
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '36.95'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '-0.55%'
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.119'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '-0.08%'
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.07807'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '0.55%'
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '8.243'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '0.28%'
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.881'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '0.06%'
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '-1.57%'
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.9256'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '0.52%'
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.1087'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '-10.08%'
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.1897'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '-0.05%'
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08903'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '-2.80%'
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.03350'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '-2.48%'
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.09588'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '-0.98%'
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.001375'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '0.75%'
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.005649'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '-5.38%'
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.449'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '-2.99%'
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '4.396'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '0.42%'
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.3421'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.258'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '18.97%'
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.1303'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '2.73%'
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.2435'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '-6.02%'
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.04325'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '-0.99%'
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.001196'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '-0.19%'
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.004247'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '-0.08%'
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0001302'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '0.09%'
$ws.Range('E26').Style = 'Normal'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '-98.10%'
$ws.Range('E27').Style = 'Normal'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.02160'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '3.16%'
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.05039'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '0.19%'
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.007504'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '-2.46%'
$ws.Range('E41').Style = 'Normal'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '0.21%'
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.008681'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '-11.27%'
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.001993'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '-3.33%'
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.008775'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '-8.31%'
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00006560'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '-2.30%'
$ws.Range('E46').Style = 'Normal'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '0.10%'
$ws.Range('E47').Style = 'Normal'
$ws.Range('B48').Value = 'CoinbaseStockToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.001002'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '-16.51%'
$ws.Range('E48').Style = 'Normal'
$ws.Range('B49').Value = 'BOLO'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.002828'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '-3.71%'
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.00002103'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '0.10%'
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0002003'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '0.10%'
$ws.Range('E51').Style = 'Normal'
